$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2137896.5
$ws.Range("I15").Value = 2137896.5
$ws.Range("K15").Value = 6413689.5
$ws.Range("M15").Value = -6413520.5
$ws.Range("H86").Value = 4716.5
$ws.Range("J86").Value = 5333.222
$ws.Range("L86").Value = 5333.222
$ws.Range("N86").Value = -7579.222
$ws.Range("H89").Value = 4716.5
$ws.Range("J89").Value = 5333.222
$ws.Range("L89").Value = 26666.11
$ws.Range("N89").Value = -37898.11
$ws.Range("H101").Value = 2978.3333
$ws.Range("I101").Value = 3664
$ws.Range("J101").Value = 2121.25
$ws.Range("K101").Value = 10992
$ws.Range("L101").Value = 6363.75
$ws.Range("M101").Value = -9370
$ws.Range("N101").Value = -9607.75
$ws.Range("H108").Value = 64293.332
$ws.Range("J108").Value = 64293.332
$ws.Range("L108").Value = 64293.332
$ws.Range("N108").Value = -71973.33199999999
$ws.Range("H111").Value = 1702
$ws.Range("I111").Value = 923.2
$ws.Range("K111").Value = 2769.6
$ws.Range("M111").Value = 297.3999999999996
$ws.Range("H112").Value = 1746.1765
$ws.Range("J112").Value = 1998.8334
$ws.Range("L112").Value = 5996.5002
$ws.Range("N112").Value = -8212.5002
$ws.Range("H135").Value = 191.66667
$ws.Range("I135").Value = 191.66667
$ws.Range("K135").Value = 1725.00003
$ws.Range("M135").Value = 809.9999699999998
$ws.Range("H137").Value = 635109.4399999999
$ws.Range("I137").Value = 1003508.2
$ws.Range("J137").Value = 21111.5
$ws.Range("K137").Value = 3010524.6
$ws.Range("L137").Value = 63334.5
$ws.Range("M137").Value = -3007974.6
$ws.Range("N137").Value = -68434.5
$ws.Range("H138").Value = 2439.9797
$ws.Range("J138").Value = 2632.7087
$ws.Range("L138").Value = 7898.1261
$ws.Range("N138").Value = -18178.1261
$ws.Range("H140").Value = 86500
$ws.Range("J140").Value = 86500
$ws.Range("L140").Value = 86500
$ws.Range("N140").Value = -96860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9959.67
$ws.Range("I32").Value = 2839.92
$ws.Range("J32").Value = 31318.92
$ws.Range("K32").Value = 2839.92
$ws.Range("L32").Value = 31318.92
$ws.Range("M32").Value = -2552.92
$ws.Range("N32").Value = -31892.92
$ws.Range("H33").Value = 18000
$ws.Range("I33").Value = 18000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 18000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -17671
$ws.Range("N33").ClearContents()
$ws.Range("H39").Value = 4624.875
$ws.Range("I39").Value = 4071.2856
$ws.Range("J39").Value = 8500
$ws.Range("K39").Value = 4071.2856
$ws.Range("L39").Value = 8500
$ws.Range("M39").Value = -3551.2856
$ws.Range("N39").Value = -9540
$ws.Range("H61").Value = 31254766
$ws.Range("I61").Value = 27781820
$ws.Range("J61").Value = 41673604
$ws.Range("K61").Value = 27781820
$ws.Range("L61").Value = 41673604
$ws.Range("M61").Value = -27781608
$ws.Range("N61").Value = -41674028
$ws.Range("H74").Value = 13901713
$ws.Range("I74").Value = 25005438
$ws.Range("K74").Value = 25005438
$ws.Range("M74").Value = -25004564
$ws.Range("H77").Value = 13901713
$ws.Range("I77").Value = 25005438
$ws.Range("K77").Value = 125027190
$ws.Range("M77").Value = -125022822
$ws.Range("H132").Value = 17551628
$ws.Range("I132").Value = 27782438
$ws.Range("J132").Value = 13097
$ws.Range("K132").Value = 83347314
$ws.Range("L132").Value = 39291
$ws.Range("M132").Value = -83344784
$ws.Range("N132").Value = -44351
$ws.Range("H136").Value = 31254766
$ws.Range("I136").Value = 27781820
$ws.Range("J136").Value = 41673604
$ws.Range("K136").Value = 83345460
$ws.Range("L136").Value = 125020812
$ws.Range("M136").Value = -83342910
$ws.Range("N136").Value = -125025912

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 46860.77
$ws.Range("J38").Value = 42497.5
$ws.Range("L38").Value = 42497.5
$ws.Range("N38").Value = -43329.5
$ws.Range("H86").Value = 3399.3157
$ws.Range("I86").Value = 3126.0908
$ws.Range("K86").Value = 3126.0908
$ws.Range("M86").Value = -2003.0908
$ws.Range("H89").Value = 3399.3157
$ws.Range("I89").Value = 3126.0908
$ws.Range("K89").Value = 15630.454
$ws.Range("M89").Value = -10014.454

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 179.4762
$ws.Range("I7").Value = 163.38461
$ws.Range("J7").Value = 205.625
$ws.Range("K7").Value = 163.38461
$ws.Range("L7").Value = 205.625
$ws.Range("M7").Value = -50.38461000000001
$ws.Range("N7").Value = -431.625
$ws.Range("H51").Value = 46663.332
$ws.Range("J51").Value = 59995
$ws.Range("L51").Value = 59995
$ws.Range("N51").Value = -61467
$ws.Range("H58").Value = 3093.7646
$ws.Range("J58").Value = 5111.5
$ws.Range("L58").Value = 5111.5
$ws.Range("N58").Value = -5517.5
$ws.Range("H61").Value = 46663.332
$ws.Range("J61").Value = 59995
$ws.Range("L61").Value = 59995
$ws.Range("N61").Value = -60691
$ws.Range("H127").Value = 92197.8
$ws.Range("J127").Value = 80247.25
$ws.Range("L127").Value = 80247.25
$ws.Range("N127").Value = -90167.25
$ws.Range("H132").Value = 5314.853
$ws.Range("I132").Value = 1752.75
$ws.Range("K132").Value = 5258.25
$ws.Range("M132").Value = -2728.25
$ws.Range("H134").Value = 2802.348
$ws.Range("I134").Value = 1654.7646
$ws.Range("K134").Value = 4964.293799999999
$ws.Range("M134").Value = -2429.293799999999
$ws.Range("H135").Value = 114998.336
$ws.Range("J135").Value = 114998.336
$ws.Range("L135").Value = 114998.336
$ws.Range("N135").Value = -125138.336
$ws.Range("H136").Value = 3093.7646
$ws.Range("J136").Value = 5111.5
$ws.Range("L136").Value = 15334.5
$ws.Range("N136").Value = -20434.5
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1100.8334
$ws.Range("I7").Value = 525.75
$ws.Range("J7").Value = 2251
$ws.Range("K7").Value = 1577.25
$ws.Range("L7").Value = 6753
$ws.Range("M7").Value = -1465.25
$ws.Range("N7").Value = -6977
$ws.Range("H40").Value = 278.72726
$ws.Range("I40").Value = 18.555555
$ws.Range("J40").Value = 1449.5
$ws.Range("K40").Value = 74.22221999999999
$ws.Range("L40").Value = 5798
$ws.Range("M40").Value = -5.222219999999993
$ws.Range("N40").Value = -5936
$ws.Range("H92").Value = 454.14285
$ws.Range("I92").Value = 410
$ws.Range("J92").Value = 487.25
$ws.Range("K92").Value = 1230
$ws.Range("L92").Value = 1461.75
$ws.Range("M92").Value = 18
$ws.Range("N92").Value = -3957.75
$ws.Range("H134").Value = 5632.75
$ws.Range("J134").Value = 14000
$ws.Range("L134").Value = 42000
$ws.Range("N134").Value = -52140
$ws.Range("H137").Value = 5273.857
$ws.Range("I137").Value = 4839.5
$ws.Range("J137").Value = 5599.625
$ws.Range("K137").Value = 14518.5
$ws.Range("L137").Value = 16798.875
$ws.Range("M137").Value = -9418.5
$ws.Range("N137").Value = -26998.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 26382818
$ws.Range("I11").Value = 23150376
$ws.Range("K11").Value = 23150376
$ws.Range("M11").Value = -23150237
$ws.Range("H55").Value = 12769.889
$ws.Range("I55").Value = 12503.75
$ws.Range("K55").Value = 12503.75
$ws.Range("M55").Value = -12176.75
$ws.Range("H70").Value = 6856
$ws.Range("I70").Value = 6489.75
$ws.Range("K70").Value = 6489.75
$ws.Range("M70").Value = -6219.75
$ws.Range("H73").Value = 6856
$ws.Range("I73").Value = 6489.75
$ws.Range("K73").Value = 6489.75
$ws.Range("M73").Value = -5553.75
$ws.Range("H80").Value = 2890.8333
$ws.Range("I80").Value = 2797
$ws.Range("J80").Value = 2926.923
$ws.Range("K80").Value = 2797
$ws.Range("L80").Value = 2926.923
$ws.Range("M80").Value = -1799
$ws.Range("N80").Value = -4922.923
$ws.Range("H83").Value = 2890.8333
$ws.Range("I83").Value = 2797
$ws.Range("J83").Value = 2926.923
$ws.Range("K83").Value = 13985
$ws.Range("L83").Value = 14634.615
$ws.Range("M83").Value = -8993
$ws.Range("N83").Value = -24618.615
$ws.Range("H132").Value = 55560136
$ws.Range("I132").Value = 55560136
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 166680408
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -166677878
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 828.03705
$ws.Range("I16").Value = 691.4
$ws.Range("J16").Value = 2536
$ws.Range("K16").Value = 691.4
$ws.Range("L16").Value = 2536
$ws.Range("M16").Value = -521.4
$ws.Range("N16").Value = -2876
$ws.Range("H55").Value = 100000770
$ws.Range("I55").Value = 125000776
$ws.Range("J55").Value = 744
$ws.Range("K55").Value = 125000776
$ws.Range("L55").Value = 744
$ws.Range("M55").Value = -125000603
$ws.Range("N55").Value = -1090
$ws.Range("H129").Value = 77195
$ws.Range("I129").Value = 66390
$ws.Range("J129").Value = 88000
$ws.Range("K129").Value = 66390
$ws.Range("L129").Value = 88000
$ws.Range("M129").Value = -61390
$ws.Range("N129").Value = -98000

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2492.25
$ws.Range("I96").Value = 2492.25
$ws.Range("K96").Value = 2492.25
$ws.Range("M96").Value = -1119.25
$ws.Range("H126").Value = 2588.5
$ws.Range("I126").Value = 2618.4
$ws.Range("J126").Value = 2558.6
$ws.Range("K126").Value = 7855.200000000001
$ws.Range("L126").Value = 7675.799999999999
$ws.Range("M126").Value = -5385.200000000001
$ws.Range("N126").Value = -12615.8
